$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# New rows to append to the bottom of the "Completed" reading list table.
# Columns: A=Title, B=Author, C=Start Date, D=Finish Date, E=Tags,
#          F=Type, G=Length, H=Rating, I=First Time Reading?
$newRows = @(
    @{ Row = 115; Title = "Salt"; Author = "Mark Kurlansky";
       Start = 44041; Finish = 44045;
       Tags = "history;salt;cuisine;food;fish";
       Type = "Audio"; Length = "13 Hours 49 Mins"; Rating = 3 },
    @{ Row = 116; Title = "Guns, Germs and Steel"; Author = "Jared Diamond";
       Start = 44045; Finish = 44048;
       Tags = "history;society;slavery;conquest;western culture;hunters & gatherers;farming;disease";
       Type = "Audio"; Length = "16 Hours 21 Mins"; Rating = 3 },
    @{ Row = 117; Title = "Get Well Soon"; Author = "Jennifer Wright";
       Start = 44048; Finish = 44049;
       Tags = "history;disease;spanish flu;pandemic;aids;cholera;polio;lobotomy;bubonic plague;dancing plague;smallpox;syphilis;tuberculosis;leprosy;encephalitis lethargica;medicine;vaccine";
       Type = "Audio"; Length = "7 Hours 49 Mins"; Rating = 4 }
)

# Mirror the order in which the data was actually typed in (Title, Author,
# dates, Tags and Type for every row, then the Length for every row, then
# Rating/First-Time-Reading) so that new shared-string entries land in the
# same sequence as the authoritative edit.
foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Title
    $ws.Cells.Item($row, 2).Value = $r.Author

    # Copy the date formatting from the existing last data row (114) so the
    # new date cells share the same style as the rest of the column.
    $ws.Cells.Item(114, 3).Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
    $ws.Cells.Item($row, 3).Value = $r.Start

    $ws.Cells.Item(114, 4).Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)
    $ws.Cells.Item($row, 4).Value = $r.Finish

    $ws.Cells.Item($row, 5).Value = $r.Tags
    $ws.Cells.Item($row, 6).Value = $r.Type
}

# The audiobook lengths were looked up and filled in out of row order
# (Guns Germs and Steel, then Get Well Soon, then Salt) -- replicate that
# exact entry order so new shared-string ids line up.
$lengthOrder = @(116, 117, 115)
foreach ($rowNum in $lengthOrder) {
    $r = $newRows | Where-Object { $_.Row -eq $rowNum }
    $ws.Cells.Item($rowNum, 7).Value = $r.Length
}

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 8).Value = $r.Rating
    $ws.Cells.Item($r.Row, 9).Value = $true
}

$excel.CutCopyMode = $false

# Update the worksheet selection to match the author's cursor position after
# appending the new rows (the view also scrolls down to keep it visible).
$wb.Windows.Item(1).ScrollRow = 91
$ws.Range("A118").Select()
